$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.929.42'
$ws.Range('E2').Value = '  +2.20%  '
$ws.Range('D3').Value = '3.187.31'
$ws.Range('E3').Value = '  +0.88%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '536.45'
$ws.Range('E5').Value = '  +1.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.82'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.532'
$ws.Range('E8').Value = '  -3.55%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.31'
$ws.Range('E9').Value = '  -0.14%  '
$ws.Range('E10').Value = '  +0.29%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.433'
$ws.Range('E11').Value = '  -1.81%  '
$ws.Range('D12').Value = '3.740.52'
$ws.Range('E12').Value = '  +1.01%  '
$ws.Range('E13').Value = '  -2.78%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.69'
$ws.Range('E14').Value = '  -0.63%  '
$ws.Range('E15').Value = '  -0.33%  '
$ws.Range('D16').Value = '59.939.54'
$ws.Range('E16').Value = '  +2.12%  '
$ws.Range('D17').Value = '3.217.80'
$ws.Range('E17').Value = '  +2.23%  '
$ws.Range('E18').Value = '  -0.75%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.23'
$ws.Range('E19').Value = '  +1.53%  '
$ws.Range('E20').Value = '  +0.80%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '368.70'
$ws.Range('E21').Value = '  -2.23%  '
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('E23').Value = '  -1.87%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '69.51'
$ws.Range('E24').Value = '  -0.45%  '
$ws.Range('E25').Value = '  +1.59%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.57'
$ws.Range('E26').Value = '  +3.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.993'
$ws.Range('E27').Value = '  -0.78%  '
$ws.Range('D28').Value = '0.0₃0872'
$ws.Range('E28').Value = '  +0.48%  '
$ws.Range('E29').Value = '  +0.38%  '
$ws.Range('E30').Value = '  +0.04%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.07'
$ws.Range('E31').Value = '  +0.56%  '
$ws.Range('E32').Value = '  +2.75%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.56'
$ws.Range('E33').Value = '  +4.62%  '
$ws.Range('E34').Value = '  +2.30%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '158.19'
$ws.Range('E35').Value = '  +0.34%  '
$ws.Range('E36').Value = '  +1.55%  '
$ws.Range('E37').Value = '  +5.36%  '
$ws.Range('D38').Value = '2.782.59'
$ws.Range('E38').Value = '  +4.94%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0709'
$ws.Range('E39').Value = '  +2.31%  '
$ws.Range('E40').Value = '  +6.55%  '
$ws.Range('E41').Value = '  -0.18%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.21'
$ws.Range('E42').Value = '  -1.84%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '39.89'
$ws.Range('E43').Value = '  +1.87%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.717'
$ws.Range('E44').Value = '  -0.66%  '
$ws.Range('E45').Value = '  +0.08%  '
$ws.Range('D46').Value = '3.231.23'
$ws.Range('E46').Value = '  +0.99%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.979'
$ws.Range('E47').Value = '  +0.06%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.14'
$ws.Range('E48').Value = '  -1.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '20.53'
$ws.Range('E49').Value = '  +2.34%  '
$ws.Range('E50').Value = '  +5.62%  '
$ws.Range('E51').Value = '  +0.02%  '
